$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(19, 2).Value = "functies CRUD in backoffice voor producten, eerste werken aan app"
$ws.Cells.Item(20, 2).Value = "Api en routes api voor communicatie met app, customers backoffice"
$ws.Cells.Item(22, 2).Value = "orders toevoegen via app, opschonen app"
$ws.Cells.Item(21, 2).Value = "inloggen op app, producten weergeven en api, database voor orders, CRUD voor customer op app"
$ws.Cells.Item(18, 2).Value = "opnieuw beginnen: databases, views, controllers, factories met faker, databases geseed."
$ws.Cells.Item(23, 2).Value = "proberen statistieken weergeven op dashboard, order overzicht klant, github repo herstarten, docs bijwerken"
$ws.Cells.Item(24, 2).Value = "Presentatie, Dossier aanvullen, poster, checklist, finale versie klaarzetten"

$ws.Cells.Item(20, 1).Value = 42594
$ws.Cells.Item(21, 1).Value = 42595
$ws.Cells.Item(22, 1).Value = 42596
$ws.Cells.Item(23, 1).Value = 42597
$ws.Cells.Item(24, 1).Value = 42598

$ws.Cells.Item(19, 3).Value = 16
$ws.Cells.Item(20, 3).Value = 14
$ws.Cells.Item(21, 3).Value = 13
$ws.Cells.Item(22, 3).Value = 12
$ws.Cells.Item(23, 3).Value = 14
$ws.Cells.Item(24, 3).Value = 13

$ws.Columns.Item(2).ColumnWidth = 87.16666666666667

$null = $ws.Range("C25").Select()

